$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates scraped from the coinranking.com refresh (crypto prices / volumes / a few
# rows where two coins swapped list position, bringing B/C text along with them).
$rowUpdates = @(
    @{ Row=2; D="27.851.11"; E="  -0.55%  " }
    @{ Row=3; D="1.909.17"; E="  +0.12%  " }
    @{ Row=4; D="1.000"; E="  -0.34%  " }
    @{ Row=5; D="312.92"; E="  -1.44%  " }
    @{ Row=6; D="0.9998"; E="  -0.35%  " }
    @{ Row=7; D="0.5003"; E="  +3.54%  " }
    @{ Row=8; D="0.3788"; E="  -0.41%  " }
    @{ Row=9; D="0.07271"; E="  -1.36%  " }
    @{ Row=10; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="21.25"; E="  +2.26%  " }
    @{ Row=11; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.9031"; E="  -3.16%  " }
    @{ Row=12; E="  -1.39%  " }
    @{ Row=13; D="1.875.89"; E="  -1.42%  " }
    @{ Row=14; D="5.471"; E="  -0.19%  " }
    @{ Row=15; D="92.26"; E="  +0.41%  " }
    @{ Row=16; D="1.001"; E="  -0.43%  " }
    @{ Row=17; D="0.000008716"; E="  -1.85%  " }
    @{ Row=18; D="1.000"; E="  -0.31%  " }
    @{ Row=19; D="27.875.72"; E="  -0.57%  " }
    @{ Row=20; E="  -0.43%  " }
    @{ Row=21; D="5.166"; E="  +0.54%  " }
    @{ Row=22; D="2.132.28"; E="  -0.59%  " }
    @{ Row=23; D="10.86"; E="  -0.48%  " }
    @{ Row=24; D="6.600"; E="  -0.59%  " }
    @{ Row=25; D="152.75"; E="  -2.02%  " }
    @{ Row=26; D="1.843"; E="  -3.82%  " }
    @{ Row=27; D="2.228"; E="  +4.87%  " }
    @{ Row=28; D="18.38"; E="  -0.59%  " }
    @{ Row=29; D="114.95"; E="  -2.03%  " }
    @{ Row=30; D="4.889"; E="  -1.90%  " }
    @{ Row=31; D="0.08975"; E="  +0.32%  " }
    @{ Row=32; D="3.178"; E="  -3.12%  " }
    @{ Row=33; D="4.813"; E="  +3.11%  " }
    @{ Row=34; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.235"; E="  -1.28%  " }
    @{ Row=35; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7857"; E="  +2.06%  " }
    @{ Row=36; D="2.660"; E="  +2.83%  " }
    @{ Row=37; E="  +1.34%  " }
    @{ Row=38; D="3.062"; E="  +2.07%  " }
    @{ Row=39; D="1.093"; E="  -1.07%  " }
    @{ Row=40; D="0.5522"; E="  +0.43%  " }
    @{ Row=41; D="0.05289"; E="  +0.31%  " }
    @{ Row=42; D="6.779"; E="  -2.06%  " }
    @{ Row=43; D="114.18"; E="  +2.91%  " }
    @{ Row=44; D="8.493"; E="  -0.02%  " }
    @{ Row=45; E="  -1.13%  " }
    @{ Row=46; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.4794"; E="  -0.49%  " }
    @{ Row=47; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="10.51"; E="  -1.90%  " }
    @{ Row=48; D="0.9995"; E="  -0.39%  " }
    @{ Row=49; D="1.635"; E="  -0.65%  " }
    @{ Row=50; D="67.15"; E="  -1.05%  " }
    @{ Row=51; D="0.06033"; E="  -0.64%  " }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Prices are stored as plain text in the sheet (inline strings), even when they
        # look numeric (e.g. '1.000' or dotted thousand-separators like '27.851.11').
        # Force text so Excel doesn't reinterpret them as numbers/dates.
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
